$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.471.88'
$ws.Range("E2").Value = '  +3.59%  '
$ws.Range("D3").Value = '1.795.64'
$ws.Range("E3").Value = '  +4.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.77'
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3791'
$ws.Range("E7").Value = '  +2.32%  '
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.42'
$ws.Range("E8").Value = '  +2.41%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3447'
$ws.Range("E9").Value = '  +2.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.212'
$ws.Range("E10").Value = '  +2.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07517'
$ws.Range("E11").Value = '  +1.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.008'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.07'
$ws.Range("E13").Value = '  +10.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.523'
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").Value = '1.787.01'
$ws.Range("E15").Value = '  +4.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.057'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001101'
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06675'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.09'
$ws.Range("E19").Value = '  +3.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.007'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.49'
$ws.Range("E21").Value = '  +6.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.514'
$ws.Range("E22").Value = '  +6.70%  '
$ws.Range("D23").Value = '27.321.61'
$ws.Range("E23").Value = '  +3.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.53'
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E25").Value = '  +1.09%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.537'
$ws.Range("E26").Value = '  +11.01%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.580'
$ws.Range("E27").Value = '  +8.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.65'
$ws.Range("E28").Value = '  +11.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.57'
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").Value = '1.985.96'
$ws.Range("E30").Value = '  +3.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.03'
$ws.Range("E31").Value = '  +2.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.094'
$ws.Range("E32").Value = '  -0.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.127'
$ws.Range("E33").Value = '  +4.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08653'
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.31'
$ws.Range("E35").Value = '  +5.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.666'
$ws.Range("E36").Value = '  -2.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.482'
$ws.Range("E37").Value = '  +2.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6926'
$ws.Range("E38").Value = '  +12.63%  '
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2217'
$ws.Range("E39").Value = '  +3.06%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06381'
$ws.Range("E40").Value = '  +3.26%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.848'
$ws.Range("E41").Value = '  +5.39%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02352'
$ws.Range("E42").Value = '  +1.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.273'
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.30'
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6496'
$ws.Range("E45").Value = '  +9.26%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.008'
$ws.Range("E46").Value = '  +0.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.858'
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.139'
$ws.Range("E48").Value = '  +5.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.68'
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07217'
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.61'
$ws.Range("E51").Value = '  +3.95%  '
